# Append new data rows (563-604) to Sheet1, mirroring the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(563, "2026/12/29", "火", 13, 9),
    @(564, "2026/12/29", "火", 16, 10),
    @(565, "2026/12/29", "火", 19, 10),
    @(566, "2026/12/29", "火", 23, 10),
    @(567, "2026/12/30", "水", 2, 10),
    @(568, "2026/12/30", "水", 5, 10),
    @(569, "2026/12/30", "水", 8, 10),
    @(570, "2026/12/30", "水", 13, 10),
    @(571, "2026/12/30", "水", 16, 9),
    @(572, "2026/12/30", "水", 22, 10),
    @(573, "2026/12/31", "木", 2, 11),
    @(574, "2026/12/31", "木", 6, 13),
    @(575, "2026/12/31", "木", 9, 13),
    @(576, "2026/12/31", "木", 12, 34),
    @(577, "2026/12/31", "木", 14, 15),
    @(578, "2026/12/31", "木", 22, 12),
    @(579, "2027/01/01", "金", 2, 13),
    @(580, "2027/01/01", "金", 5, 12),
    @(581, "2027/01/01", "金", 13, 14),
    @(582, "2027/01/01", "金", 16, 11),
    @(583, "2027/01/01", "金", 19, 13),
    @(584, "2027/01/02", "土", 1, 12),
    @(585, "2027/01/02", "土", 5, 12),
    @(586, "2027/01/02", "土", 8, 13),
    @(587, "2027/01/02", "土", 13, 16),
    @(588, "2027/01/02", "土", 16, 19),
    @(589, "2027/01/02", "土", 19, 21),
    @(590, "2027/01/02", "土", 22, 22),
    @(591, "2027/01/03", "日", 1, 23),
    @(592, "2027/01/03", "日", 4, 26),
    @(593, "2027/01/03", "日", 7, 23),
    @(594, "2027/01/03", "日", 13, 23),
    @(595, "2027/01/03", "日", 16, 24),
    @(596, "2027/01/03", "日", 19, 26),
    @(597, "2027/01/03", "日", 22, 21),
    @(598, "2027/01/04", "月", 2, 19),
    @(599, "2027/01/04", "月", 4, 18),
    @(600, "2027/01/04", "月", 7, 19),
    @(601, "2027/01/04", "月", 13, 20),
    @(602, "2027/01/04", "月", 22, 13),
    @(603, "2027/01/05", "火", 1, 13),
    @(604, "2027/01/05", "火", 7, 14)
)

foreach ($entry in $newRows) {
    $r = $entry[0]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $entry[1]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}
